$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 10)
$ws.Range("A10").Value = "Efrain"
$ws.Range("B10").Value = "Pinto Mancera"
$ws.Range("C10").Value = "Alfredo De La Peña, Cra 7 Cl 15#50"
$ws.Range("D10").Value = "efrainpintomancera@gmail.com"
$ws.Range("E10").Value = 8541524514

# Apply hyperlink (same mailto pattern as D8/D9), then reapply the hyperlink cell style
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:efrainpintomancera@gmail.com") | Out-Null
$ws.Range("D10").Style = $ws.Range("D9").Style

# Update the view: scroll so column B is the top-left visible, select E9
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E9").Select() | Out-Null

$wb.Save()
